$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row cells: "<name>_old" -> "<name>_FV2410", "<name>_new" -> "<name>_FV2504"
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Text
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2410")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2504")
    }
}

# Freeze the header row (pane split after row 1, active pane bottom-left)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a table ("Table1") so the header row is recognised as column headers
$rng = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
